$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "66.099.78"
$ws.Range("E2").Value = "  +4.07%  "

# Row 3
$ws.Range("D3").Value = "3.796.86"
$ws.Range("E3").Value = "  +6.72%  "

# Row 4
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
Set-TextValue $ws "D5" "428.64"
$ws.Range("E5").Value = "  +8.07%  "

# Row 6
Set-TextValue $ws "D6" "138.76"
$ws.Range("E6").Value = "  +12.79%  "

# Row 7
$ws.Range("E7").Value = "  +6.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
Set-TextValue $ws "D9" "0.741"
$ws.Range("E9").Value = "  +8.91%  "

# Row 10
$ws.Range("E10").Value = "  +1.07%  "

# Row 11
Set-TextValue $ws "D11" "0.0000314"
$ws.Range("E11").Value = "  -3.44%  "

# Row 12
Set-TextValue $ws "D12" "42.94"

# Row 13
Set-TextValue $ws "D13" "10.55"
$ws.Range("E13").Value = "  +14.87%  "

# Row 14
$ws.Range("D14").Value = "4.387.37"
$ws.Range("E14").Value = "  +6.44%  "

# Row 15
Set-TextValue $ws "D15" "15.06"
$ws.Range("E15").Value = "  +8.99%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D16" "0.138"
$ws.Range("E16").Value = "  +1.21%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.801.59"
$ws.Range("E17").Value = "  +7.45%  "

# Row 18
Set-TextValue $ws "D18" "20.02"
$ws.Range("E18").Value = "  +7.10%  "

# Row 19
$ws.Range("E19").Value = "  +10.91%  "

# Row 20
$ws.Range("D20").Value = "66.216.10"
$ws.Range("E20").Value = "  +4.12%  "

# Row 21
Set-TextValue $ws "D21" "407.12"
$ws.Range("E21").Value = "  +3.32%  "

# Row 22
Set-TextValue $ws "D22" "15.18"
$ws.Range("E22").Value = "  +9.27%  "

# Row 23
$ws.Range("E23").Value = "  +10.90%  "

# Row 24
Set-TextValue $ws "D24" "85.25"
$ws.Range("E24").Value = "  +4.04%  "

# Row 25
Set-TextValue $ws "D25" "36.78"
$ws.Range("E25").Value = "  +8.33%  "

# Row 26
Set-TextValue $ws "D26" "9.98"
$ws.Range("E26").Value = "  +46.98%  "

# Row 27
Set-TextValue $ws "D27" "3.30"

# Row 28
$ws.Range("E28").Value = "  +13.40%  "

# Row 29
Set-TextValue $ws "D29" "5.42"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
Set-TextValue $ws "D30" "13.90"
$ws.Range("E30").Value = "  +16.01%  "

# Row 31
Set-TextValue $ws "D31" "710.86"
$ws.Range("E31").Value = "  +4.32%  "

# Row 32
Set-TextValue $ws "D32" "0.132"
$ws.Range("E32").Value = "  +17.94%  "

# Row 33
$ws.Range("E33").Value = "  +7.46%  "

# Row 34
Set-TextValue $ws "D34" "41.15"
$ws.Range("E34").Value = "  +12.34%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
Set-TextValue $ws "D36" "5.72"
$ws.Range("E36").Value = "  +40.95%  "

# Row 37
Set-TextValue $ws "D37" "0.150"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
Set-TextValue $ws "D38" "56.66"
$ws.Range("E38").Value = "  +5.14%  "

# Row 39
Set-TextValue $ws "D39" "0.0477"
$ws.Range("E39").Value = "  +9.26%  "

# Row 40
$ws.Range("E40").Value = "  +48.89%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0685"
$ws.Range("E41").Value = "  +4.92%  "

# Row 42
$ws.Range("E42").Value = "  +8.70%  "

# Row 43
$ws.Range("E43").Value = "  +6.83%  "

# Row 44
$ws.Range("E44").Value = "  +0.19%  "

# Row 45
Set-TextValue $ws "D45" "3.39"
$ws.Range("E45").Value = "  +10.33%  "

# Row 46
Set-TextValue $ws "D46" "0.323"
$ws.Range("E46").Value = "  +16.75%  "

# Row 47
Set-TextValue $ws "D47" "3.13"
$ws.Range("E47").Value = "  +1.73%  "

# Row 48
$ws.Range("E48").Value = "  +6.73%  "

# Row 49
$ws.Range("E49").Value = "  +6.01%  "

# Row 50
Set-TextValue $ws "D50" "142.70"
$ws.Range("E50").Value = "  +0.68%  "

# Row 51
$ws.Range("E51").Value = "  +5.98%  "
